# Add dedicated DATASETS threshold columns (threshold_method, threshold_min_size,
# threshold_pct, threshold_z, threshold_window_days) and push the existing
# "notes" column from M out to R, moving its data along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATASETS")

# --- Header row (row 1) -----------------------------------------------------
# Preserve the bold "header" formatting that currently lives on M1 by copying
# it onto the new notes column (R1) before M1's text/format changes.
$ws.Range("M1").Copy()
$ws.Range("R1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R1").Value = "notes"

# M1 becomes a plain (non-bold) header cell for the new threshold_method column.
$ws.Range("M1").ClearFormats()
$ws.Range("M1").Value = "threshold_method"

$ws.Range("N1").Value = "threshold_min_size"
$ws.Range("O1").Value = "threshold_pct"
$ws.Range("P1").Value = "threshold_z"
$ws.Range("Q1").Value = "threshold_window_days"

# --- Plain data rows: move the free-text notes from M to R ------------------
$ws.Range("M2").Cut($ws.Range("R2"))
$ws.Range("M3").Cut($ws.Range("R3"))
$ws.Range("M4").Cut($ws.Range("R4"))
$ws.Range("M5").Cut($ws.Range("R5"))
$ws.Range("M6").Cut($ws.Range("R6"))
$ws.Range("M8").Cut($ws.Range("R8"))
$ws.Range("M10").Cut($ws.Range("R10"))
$ws.Range("M12").Cut($ws.Range("R12"))

# --- ES_BIG_TRADES (row 14): split the old free-text threshold info into the
# new dedicated columns, leaving only the instrument_id reference in notes. ---
$ws.Range("R14").Value = "instrument_id: ES"
$ws.Range("M14").Value = "fixed_count"
$ws.Range("N14").Value = 50

# --- ES_BIG_TRADES_PROXY (row 16): same treatment. --------------------------
$ws.Range("R16").Value = "instrument_id: ES"
$ws.Range("M16").Value = "fixed_count"
$ws.Range("N16").Value = 100

# --- Extend the sheet's used range down to row 17 (matches the committed
# dimension A1:R17) without leaving any visible content behind. -------------
$ws.Range("A1").Copy()
$ws.Range("R17").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R17").ClearFormats()
